$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row gets a new 4th column "posicao" ---
$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "inicial"
$ws.Range("D1").Value = "posicao"

# --- Data rows: HANTAROGAMER/RKFox rows swap places (now sorted by
#     "inicial" descending) and each row gets its new "posicao" rank ---
$ws.Range("A2").Value = "61e484ca5aa1be001868f065"
$ws.Range("B2").Value = "HANTAROGAMER"
$ws.Range("C2").Value = 1800000000
$ws.Range("D2").Value = 1

$ws.Range("A3").Value = "61e852b4dc27dc001969efa3"
$ws.Range("B3").Value = "RKFox"
$ws.Range("C3").Value = 800000000
$ws.Range("D3").Value = 2

$ws.Range("A4").Value = "65de9e82a1e9f41193e2f6cc"
$ws.Range("B4").Value = "GUERDE"
$ws.Range("C4").Value = 80000000
$ws.Range("D4").Value = 3

# --- Sort the data (descending by "inicial", column C) ---
$sortRange = $ws.Range("A1:D4")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("C1"), 0, 2, 0, 0)
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 1
$ws.Sort.Apply()

# --- Turn on AutoFilter over the header row ---
$ws.Range("A1:D1").AutoFilter()

# Excel normally registers the hidden _FilterDatabase defined name
# whenever AutoFilter is switched on for a sheet.
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=Planilha1!`$A`$1:`$D`$1")
$filterName.Visible = $false

# --- Selection moves to D5, just below the new column ---
$ws.Range("D5").Select()
